$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
for ($r = 1; $r -le 31; $r++) {
    if ($r -eq 14) {
        $ws.Cells.Item($r, 2).Value2 = $ws.Cells.Item($r, 4).Value2
        $ws.Cells.Item($r, 4).Value2 = $null
    } else {
        $ws.Cells.Item($r, 2).Value2 = $ws.Cells.Item($r, 1).Value2
    }
    $ws.Rows.Item($r).AutoFit()
}
